$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: clear the previously-blank cells E5, F5, G5, L5 so they drop out of the sheet entirely
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = ""
$ws.Range("L5").Value = ""

# New row 6: duplicate of the "Test Ringover (NO TOCAR)" record with only Numero/Nombre/Fecha de venta filled in
$ws.Range("A6").Value = 2488
$ws.Range("B6").Value = "Test Ringover (NO TOCAR)"
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = ""
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("L6").Value = ""
$ws.Range("M6").Value = ""
$ws.Range("N6").Value = "2024-01-03T10:49:29.104Z"
